$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append: date, health_reg_name, admissions
$newRows = @(
    @(43927, "Helse Midt-Norge", 6),
    @(43927, "Helse Nord", 6),
    @(43927, "Helse Sør-Øst", 60),
    @(43927, "Helse Vest", 11)
)

$startRow = 98
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Value = $rowData[0]
    $dateCell.NumberFormat = "yyyy-mm-dd"

    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
}
